$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 6111.2964
$ws.Range("I15").Value = 6111.2964
$ws.Range("K15").Value = 18333.8892
$ws.Range("M15").Value = -18164.8892
$ws.Range("H86").Value = 14464
$ws.Range("I86").Value = 16124.667
$ws.Range("K86").Value = 16124.667
$ws.Range("M86").Value = -15001.667
$ws.Range("H89").Value = 14464
$ws.Range("I89").Value = 16124.667
$ws.Range("K89").Value = 80623.33499999999
$ws.Range("M89").Value = -75007.33499999999
$ws.Range("H107").Value = 443.3125
$ws.Range("I107").Value = 446.07693
$ws.Range("J107").Value = 431.33334
$ws.Range("K107").Value = 446.07693
$ws.Range("L107").Value = 431.33334
$ws.Range("M107").Value = 1473.92307
$ws.Range("N107").Value = -4271.33334
$ws.Range("H113").Value = 58836310
$ws.Range("I113").Value = 71432450
$ws.Range("K113").Value = 71432450
$ws.Range("M113").Value = -71429196
$ws.Range("H137").Value = 2212.8
$ws.Range("I137").Value = 2022.7273
$ws.Range("K137").Value = 6068.1819
$ws.Range("M137").Value = -3518.1819
$ws.Range("H138").Value = 3547.2
$ws.Range("J138").Value = 3167.625
$ws.Range("L138").Value = 9502.875
$ws.Range("N138").Value = -19782.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6467.6055
$ws.Range("I32").Value = 6639.108
$ws.Range("J32").Value = 122
$ws.Range("K32").Value = 6639.108
$ws.Range("L32").Value = 122
$ws.Range("M32").Value = -6352.108
$ws.Range("N32").Value = -696
$ws.Range("H97").Value = 1757.7916
$ws.Range("I97").Value = 1636.15
$ws.Range("K97").Value = 1636.15
$ws.Range("M97").Value = -1140.15
$ws.Range("H110").Value = 2667.4666
$ws.Range("I110").Value = 918.5833
$ws.Range("J110").Value = 9663
$ws.Range("K110").Value = 918.5833
$ws.Range("L110").Value = 9663
$ws.Range("M110").Value = 1126.4167
$ws.Range("N110").Value = -13753

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1802.3334
$ws.Range("J86").Value = 1841.6
$ws.Range("L86").Value = 1841.6
$ws.Range("N86").Value = -4087.6
$ws.Range("H89").Value = 1802.3334
$ws.Range("J89").Value = 1841.6
$ws.Range("L89").Value = 9208
$ws.Range("N89").Value = -20440
$ws.Range("H134").Value = 3234.3784
$ws.Range("I134").Value = 2848.0857
$ws.Range("J134").Value = 9994.5
$ws.Range("K134").Value = 8544.257100000001
$ws.Range("L134").Value = 29983.5
$ws.Range("M134").Value = -6009.257100000001
$ws.Range("N134").Value = -35053.5
$ws.Range("H140").Value = 85662.664
$ws.Range("J140").Value = 85662.664
$ws.Range("L140").Value = 85662.664
$ws.Range("N140").Value = -96022.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2276.5
$ws.Range("J5").Value = 3364.75
$ws.Range("L5").Value = 3364.75
$ws.Range("N5").Value = -3588.75
$ws.Range("H107").Value = 755.65216
$ws.Range("I107").Value = 331.66666
$ws.Range("K107").Value = 331.66666
$ws.Range("M107").Value = 1588.33334
$ws.Range("H141").Value = 138265
$ws.Range("J141").Value = 152205.19
$ws.Range("L141").Value = 152205.19
$ws.Range("N141").Value = -162565.19

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 924.3077
$ws.Range("I50").Value = 1671.75
$ws.Range("J50").Value = 592.1111
$ws.Range("K50").Value = 5015.25
$ws.Range("L50").Value = 1776.3333
$ws.Range("M50").Value = -4534.25
$ws.Range("N50").Value = -2738.3333
$ws.Range("H53").Value = 924.3077
$ws.Range("I53").Value = 1671.75
$ws.Range("J53").Value = 592.1111
$ws.Range("K53").Value = 5015.25
$ws.Range("L53").Value = 1776.3333
$ws.Range("M53").Value = -4534.25
$ws.Range("N53").Value = -2738.3333
$ws.Range("H131").Value = 31253478
$ws.Range("J131").Value = 3954.5833
$ws.Range("L131").Value = 11863.7499
$ws.Range("N131").Value = -21943.7499

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H80").Value = 2742.2856
$ws.Range("J80").Value = 2965.6667
$ws.Range("L80").Value = 2965.6667
$ws.Range("N80").Value = -4961.6667
$ws.Range("H83").Value = 2742.2856
$ws.Range("J83").Value = 2965.6667
$ws.Range("L83").Value = 14828.3335
$ws.Range("N83").Value = -24812.3335
$ws.Range("H126").Value = 4831.1665
$ws.Range("J126").Value = 9999.666999999999
$ws.Range("L126").Value = 29999.001
$ws.Range("N126").Value = -34939.001
$ws.Range("H132").Value = 2512.9285
$ws.Range("I132").Value = 2512.9285
$ws.Range("K132").Value = 7538.7855
$ws.Range("M132").Value = -5008.7855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2432.6667
$ws.Range("J10").Value = 2432.6667
$ws.Range("L10").Value = 2432.6667
$ws.Range("N10").Value = -2712.6667
$ws.Range("H22").Value = 2628.1428
$ws.Range("I22").Value = 2399
$ws.Range("J22").Value = 2800
$ws.Range("K22").Value = 2399
$ws.Range("L22").Value = 2800
$ws.Range("M22").Value = -2104
$ws.Range("N22").Value = -3390
$ws.Range("H27").Value = 2628.1428
$ws.Range("I27").Value = 2399
$ws.Range("J27").Value = 2800
$ws.Range("K27").Value = 2399
$ws.Range("L27").Value = 2800
$ws.Range("M27").Value = -2292
$ws.Range("N27").Value = -3014
$ws.Range("H32").Value = 12000
$ws.Range("I32").Value = 12000
$ws.Range("K32").Value = 12000
$ws.Range("M32").Value = -11683
$ws.Range("H123").Value = 55450
$ws.Range("J123").Value = 55450
$ws.Range("L123").Value = 55450
$ws.Range("N123").Value = -65250
$ws.Range("H136").Value = 95241510
$ws.Range("I136").Value = 62502924
$ws.Range("K136").Value = 187508772
$ws.Range("M136").Value = -187506222

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -860
$ws.Range("H46").Value = 57895.855
$ws.Range("J46").Value = 57895.855
$ws.Range("L46").Value = 57895.855
$ws.Range("N46").Value = -58357.855
$ws.Range("H134").Value = 57895.855
$ws.Range("J134").Value = 57895.855
$ws.Range("L134").Value = 173687.565
$ws.Range("N134").Value = -178757.565
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200
